# Update the division-problem answers in the single table, cell by cell.
# Word COM uses 1-based indices for Rows/Columns in Table.Cell(row, col).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col) -> new text, using 1-based Word indices.
# Data rows in the document are rows 1, 5, 9, 13, 17 (0-based 0,4,8,12,16).
$updates = @(
    @{ Row = 1;  Col = 1; Text = "147÷7=21, 0" },
    @{ Row = 1;  Col = 2; Text = "209÷7=29, 6" },
    @{ Row = 1;  Col = 3; Text = "188÷5=37, 3" },
    @{ Row = 1;  Col = 4; Text = "886÷3=295, 1" },
    @{ Row = 1;  Col = 5; Text = "794÷3=264, 2" },

    @{ Row = 5;  Col = 1; Text = "982÷7=140, 2" },
    @{ Row = 5;  Col = 2; Text = "709÷8=88, 5" },
    @{ Row = 5;  Col = 3; Text = "291÷4=72, 3" },
    @{ Row = 5;  Col = 4; Text = "230÷2=115, 0" },
    @{ Row = 5;  Col = 5; Text = "734÷6=122, 2" },

    @{ Row = 9;  Col = 1; Text = "871÷8=108, 7" },
    @{ Row = 9;  Col = 2; Text = "910÷4=227, 2" },
    @{ Row = 9;  Col = 3; Text = "357÷5=71, 2" },
    @{ Row = 9;  Col = 4; Text = "187÷6=31, 1" },
    @{ Row = 9;  Col = 5; Text = "232÷6=38, 4" },

    @{ Row = 13; Col = 1; Text = "924÷9=102, 6" },
    @{ Row = 13; Col = 2; Text = "939÷7=134, 1" },
    @{ Row = 13; Col = 3; Text = "872÷8=109, 0" },
    @{ Row = 13; Col = 4; Text = "588÷3=196, 0" },
    @{ Row = 13; Col = 5; Text = "712÷6=118, 4" },

    @{ Row = 17; Col = 1; Text = "241÷7=34, 3" },
    @{ Row = 17; Col = 2; Text = "741÷3=247, 0" },
    @{ Row = 17; Col = 3; Text = "403÷9=44, 7" },
    @{ Row = 17; Col = 4; Text = "488÷8=61, 0" },
    @{ Row = 17; Col = 5; Text = "929÷3=309, 2" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
